$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple numeric/percent updates (Column D = Price, Column E = Volume(1h))
$updates = @(
    @{ Row = 2;  D = "68.033.61"; E = "  +1.55%  " },
    @{ Row = 3;  D = "3.340.77";  E = "  +1.73%  " },
    @{ Row = 4;  D = "1.00";      E = "  +0.04%  " },
    @{ Row = 5;  D = "581.61";    E = "  +1.81%  " },
    @{ Row = 6;  D = "177.28";    E = "  +0.68%  " },
    @{ Row = 7;  E = "  -0.05%  " },
    @{ Row = 8;  D = "0.590";     E = "  +1.54%  " },
    @{ Row = 9;  D = "3.337.58";  E = "  +1.83%  " },
    @{ Row = 10; D = "0.182";     E = "  +5.04%  " },
    @{ Row = 11; E = "  +1.52%  " },
    @{ Row = 12; D = "46.87";     E = "  +2.42%  " },
    @{ Row = 13; D = "0.0000273"; E = "  +1.70%  " },
    @{ Row = 14; D = "690.85";    E = "  -0.59%  " },
    @{ Row = 15; D = "3.883.81";  E = "  +1.91%  " },
    @{ Row = 16; D = "8.45" },
    @{ Row = 17; D = "68.040.44"; E = "  +1.43%  " },
    @{ Row = 18; E = "  -0.32%  " },
    @{ Row = 19; D = "3.344.00";  E = "  +1.76%  " },
    @{ Row = 20; E = "  +0.61%  " },
    @{ Row = 21; D = "11.09";     E = "  +3.42%  " },
    @{ Row = 22; D = "0.898";     E = "  +1.17%  " },
    @{ Row = 23; D = "5.35";      E = "  +4.66%  " },
    @{ Row = 24; D = "17.01";     E = "  +0.27%  " },
    @{ Row = 25; D = "99.19";     E = "  -0.03%  " },
    @{ Row = 26; E = "  +0.28%  " },
    @{ Row = 27; E = "  -0.17%  " },
    @{ Row = 28; E = "  +2.62%  " },
    @{ Row = 29; D = "32.98";     E = "  -0.17%  " },
    @{ Row = 30; E = "  +2.18%  " },
    @{ Row = 31; D = "7.07";      E = "  +4.88%  " },
    @{ Row = 32; D = "572.63";    E = "  +0.31%  " },
    @{ Row = 33; E = "  +1.82%  " },
    @{ Row = 34; E = "  +2.17%  " },
    @{ Row = 38; E = "  -1.35%  " },
    @{ Row = 39; D = "34.79";     E = "  +9.22%  " },
    @{ Row = 40; E = "  +1.94%  " },
    @{ Row = 41; D = "2.66";      E = "  +1.88%  " },
    @{ Row = 42; E = "  +5.81%  " },
    @{ Row = 43; E = "  +0.60%  " },
    @{ Row = 46; E = "  +0.36%  " },
    @{ Row = 47; D = "2.66";      E = "  +5.23%  " },
    @{ Row = 48; E = "  +1.14%  " },
    @{ Row = 49; E = "  -0.30%  " },
    @{ Row = 50; E = "  -2.95%  " },
    @{ Row = 51; D = "129.84";    E = "  -0.41%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("D")) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}

# Rows 35-37 rotate: OKB/Dai/Maker -> Maker/OKB/Dai (with new price/volume values)
$ws.Cells.Item(35, 2).Value = "Maker"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.721.30"
$ws.Cells.Item(35, 5).Value = "  -4.42%  "

$ws.Cells.Item(36, 2).Value = "OKB"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "57.32"
$ws.Cells.Item(36, 5).Value = "  +3.16%  "

$ws.Cells.Item(37, 2).Value = "Dai"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.00"
$ws.Cells.Item(37, 5).Value = "  +0.06%  "

# Rows 44-45 swap: ApeXProtocol/TheGraph -> TheGraph/ApeXProtocol (with new price/volume values)
$ws.Cells.Item(44, 2).Value = "TheGraph"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.337"
$ws.Cells.Item(44, 5).Value = "  +2.66%  "

$ws.Cells.Item(45, 2).Value = "ApeXProtocol"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.31"
$ws.Cells.Item(45, 5).Value = "  -1.40%  "
